# Generate Report for Handoff
#
# Adds two new localization entries -
#   955745b8-040b-4a9f-9646-605fe0c61001.md  (inserted as the new 2nd data row)
#   ebe75446-2550-4555-a917-027670c0007d.md  (appended as the new last data row)
# - to all three worksheets (Overview, zh-cn, de-de) and their backing tables.

$wb = $excel.ActiveWorkbook

# ======================================================================
# Sheet 1: "Overview"
# ======================================================================
$ws1 = $wb.Worksheets.Item(1)

# Insert row 3 (pushes the existing "a3b6c57b..." row down to row 4, and
# inherits its cell formatting/styles for the new blank row).
$ws1.Rows.Item(3).Insert()
# Insert row 5 right after the (now) last data row 4, inheriting its format.
$ws1.Rows.Item(5).Insert()

# New row 3: 955745b8-...
$ws1.Cells.Item(3,1).Value = '955745b8-040b-4a9f-9646-605fe0c61001.md'
$ws1.Cells.Item(3,2).Value = 'e2e\955745b8-040b-4a9f-9646-605fe0c61001.md'
$ws1.Cells.Item(3,3).Value = '.md'
$ws1.Cells.Item(3,5).Value = 'Ready for handoff'
$ws1.Cells.Item(3,6).Value = 'Ready for handoff'
$ws1.Cells.Item(3,7).Value = '2016-09-06 18:54:37'

# New row 5: ebe75446-...
$ws1.Cells.Item(5,1).Value = 'ebe75446-2550-4555-a917-027670c0007d.md'
$ws1.Cells.Item(5,2).Value = 'e2e\ebe75446-2550-4555-a917-027670c0007d.md'
$ws1.Cells.Item(5,3).Value = '.md'
$ws1.Cells.Item(5,5).Value = 'Ready for handoff'
$ws1.Cells.Item(5,6).Value = 'Ready for handoff'
$ws1.Cells.Item(5,7).Value = '2016-09-06 18:54:37'

# Grow the "Overview" table to cover the two new rows.
$lo3 = $ws1.ListObjects.Item(1)
$lo3.Resize($ws1.Range("A1:G5"))

# Hyperlinks in column B shifted around by the inserts above, so rebuild
# them all against their final addresses.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Cells.Item(2,2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96b0aa7c191bf3b4f4c8cb6886b752d53e9713ae/e2e/f2388c28-632b-4c28-9359-b42de4a9fbdc.md", "", "", "e2e\f2388c28-632b-4c28-9359-b42de4a9fbdc.md")
$ws1.Hyperlinks.Add($ws1.Cells.Item(3,2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/955745b8040b4a9f9646605fe0c61001/e2e/955745b8-040b-4a9f-9646-605fe0c61001.md", "", "", "e2e\955745b8-040b-4a9f-9646-605fe0c61001.md")
$ws1.Hyperlinks.Add($ws1.Cells.Item(4,2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b504e4d7ee380ccc21786bc208918f6a2dd49eb/e2e/a3b6c57b-a213-4974-8a52-4673e4d3be66.md", "", "", "e2e\a3b6c57b-a213-4974-8a52-4673e4d3be66.md")
$ws1.Hyperlinks.Add($ws1.Cells.Item(5,2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ebe754462550455 5a917027670c0007d/e2e/ebe75446-2550-4555-a917-027670c0007d.md", "", "", "e2e\ebe75446-2550-4555-a917-027670c0007d.md")

# ======================================================================
# Sheet 2: "zh-cn"
# ======================================================================
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(3).Insert()
$ws2.Rows.Item(5).Insert()

# New row 3: 955745b8-...
$ws2.Cells.Item(3,1).Value = '955745b8-040b-4a9f-9646-605fe0c61001.md'
$ws2.Cells.Item(3,2).Value = '.md'
$ws2.Cells.Item(3,3).Value = 'Ready for handoff'
$ws2.Cells.Item(3,4).Value = 'e2e'
$ws2.Cells.Item(3,5).Value = 'ht'
$ws2.Cells.Item(3,6).Value = "'False"
$ws2.Cells.Item(3,7).Value = '955745b8-040b-4a9f-9646-605fe0c61001.fd6e2542777b073d9ea9686ccef660a638ae5ab9.zh-cn.xlf'
$ws2.Cells.Item(3,8).Value = '2016-09-06 18:54:32'
$ws2.Cells.Item(3,11).Value = '0001-01-01 00:00:00'
$ws2.Cells.Item(3,13).Value = "'True"
$ws2.Cells.Item(3,15).Value = "'False"

# New row 5: ebe75446-...
$ws2.Cells.Item(5,1).Value = 'ebe75446-2550-4555-a917-027670c0007d.md'
$ws2.Cells.Item(5,2).Value = '.md'
$ws2.Cells.Item(5,3).Value = 'Ready for handoff'
$ws2.Cells.Item(5,4).Value = 'e2e'
$ws2.Cells.Item(5,5).Value = 'ht'
$ws2.Cells.Item(5,6).Value = "'False"
$ws2.Cells.Item(5,7).Value = 'ebe75446-2550-4555-a917-027670c0007d.365d0372f7bf2916439d03b9b4bb7b8f032154a2.zh-cn.xlf'
$ws2.Cells.Item(5,8).Value = '2016-09-06 18:54:32'
$ws2.Cells.Item(5,11).Value = '0001-01-01 00:00:00'
$ws2.Cells.Item(5,13).Value = "'True"
$ws2.Cells.Item(5,15).Value = "'False"

$lo1 = $ws2.ListObjects.Item(1)
$lo1.Resize($ws2.Range("A1:P5"))

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96b0aa7c191bf3b4f4c8cb6886b752d53e9713ae/e2e/f2388c28-632b-4c28-9359-b42de4a9fbdc.md", "", "", "f2388c28-632b-4c28-9359-b42de4a9fbdc.md")
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,9), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d2f7c5922bd77470b794f6697bdc047e49ea96bf/e2e/f2388c28-632b-4c28-9359-b42de4a9fbdc.md", "", "", "f2388c28-632b-4c28-9359-b42de4a9fbdc.md")
$ws2.Hyperlinks.Add($ws2.Cells.Item(3,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/955745b8040b4a9f9646605fe0c61001/e2e/955745b8-040b-4a9f-9646-605fe0c61001.md", "", "", "955745b8-040b-4a9f-9646-605fe0c61001.md")
$ws2.Hyperlinks.Add($ws2.Cells.Item(4,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b504e4d7ee380ccc21786bc208918f6a2dd49eb/e2e/a3b6c57b-a213-4974-8a52-4673e4d3be66.md", "", "", "a3b6c57b-a213-4974-8a52-4673e4d3be66.md")
$ws2.Hyperlinks.Add($ws2.Cells.Item(5,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ebe754462550455 5a917027670c0007d/e2e/ebe75446-2550-4555-a917-027670c0007d.md", "", "", "ebe75446-2550-4555-a917-027670c0007d.md")

# ======================================================================
# Sheet 3: "de-de"
# ======================================================================
$ws3 = $wb.Worksheets.Item(3)

$ws3.Rows.Item(3).Insert()
$ws3.Rows.Item(5).Insert()

# New row 3: 955745b8-...
$ws3.Cells.Item(3,1).Value = '955745b8-040b-4a9f-9646-605fe0c61001.md'
$ws3.Cells.Item(3,2).Value = '.md'
$ws3.Cells.Item(3,3).Value = 'Ready for handoff'
$ws3.Cells.Item(3,4).Value = 'e2e'
$ws3.Cells.Item(3,5).Value = 'ht'
$ws3.Cells.Item(3,6).Value = "'False"
$ws3.Cells.Item(3,7).Value = '955745b8-040b-4a9f-9646-605fe0c61001.fd6e2542777b073d9ea9686ccef660a638ae5ab9.de-de.xlf'
$ws3.Cells.Item(3,8).Value = '2016-09-06 18:54:37'
$ws3.Cells.Item(3,11).Value = '0001-01-01 00:00:00'
$ws3.Cells.Item(3,13).Value = "'True"
$ws3.Cells.Item(3,15).Value = "'False"

# New row 5: ebe75446-...
$ws3.Cells.Item(5,1).Value = 'ebe75446-2550-4555-a917-027670c0007d.md'
$ws3.Cells.Item(5,2).Value = '.md'
$ws3.Cells.Item(5,3).Value = 'Ready for handoff'
$ws3.Cells.Item(5,4).Value = 'e2e'
$ws3.Cells.Item(5,5).Value = 'ht'
$ws3.Cells.Item(5,6).Value = "'False"
$ws3.Cells.Item(5,7).Value = 'ebe75446-2550-4555-a917-027670c0007d.365d0372f7bf2916439d03b9b4bb7b8f032154a2.de-de.xlf'
$ws3.Cells.Item(5,8).Value = '2016-09-06 18:54:37'
$ws3.Cells.Item(5,11).Value = '0001-01-01 00:00:00'
$ws3.Cells.Item(5,13).Value = "'True"
$ws3.Cells.Item(5,15).Value = "'False"

$lo2 = $ws3.ListObjects.Item(1)
$lo2.Resize($ws3.Range("A1:P5"))

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96b0aa7c191bf3b4f4c8cb6886b752d53e9713ae/e2e/f2388c28-632b-4c28-9359-b42de4a9fbdc.md", "", "", "f2388c28-632b-4c28-9359-b42de4a9fbdc.md")
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,9), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a701e45bb92c88875c82dfd416bed04f9708fa47/e2e/f2388c28-632b-4c28-9359-b42de4a9fbdc.md", "", "", "f2388c28-632b-4c28-9359-b42de4a9fbdc.md")
$ws3.Hyperlinks.Add($ws3.Cells.Item(3,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/955745b8040b4a9f9646605fe0c61001/e2e/955745b8-040b-4a9f-9646-605fe0c61001.md", "", "", "955745b8-040b-4a9f-9646-605fe0c61001.md")
$ws3.Hyperlinks.Add($ws3.Cells.Item(4,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b504e4d7ee380ccc21786bc208918f6a2dd49eb/e2e/a3b6c57b-a213-4974-8a52-4673e4d3be66.md", "", "", "a3b6c57b-a213-4974-8a52-4673e4d3be66.md")
$ws3.Hyperlinks.Add($ws3.Cells.Item(5,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ebe754462550455 5a917027670c0007d/e2e/ebe75446-2550-4555-a917-027670c0007d.md", "", "", "ebe75446-2550-4555-a917-027670c0007d.md")
